$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(21).Insert()
$ws.Range("B21").Value = "Person Index ID"
$ws.Range("C21").Value = "A unique index identifier for a person"
$ws.Rows.Item(21).RowHeight = 30
$ws.Range("D21").Select() | Out-Null
